$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Italicise the three "Assess the inherent risk" sub-bullets.
# ---------------------------------------------------------------------
$italicTexts = @(
    "Evaluating the possible consequence of the hazard using the Consequence Scale (see below)",
    "Evaluate the likelihood of that consequence using the Likelihood Scale (see below)",
    "Determine risk rating of each hazard using the Risk Matrix."
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    foreach ($needle in $italicTexts) {
        if ($t -eq $needle) {
            $rng = $p.Range
            $rng.Font.Italic = 1
            $rng.Font.ItalicBi = 1
        }
    }
}

# ---------------------------------------------------------------------
# 2. Remove the empty trailing sub-bullet paragraph that used to follow
#    "Determine risk rating of each hazard using the Risk Matrix."
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Determine risk rating of each hazard using the Risk Matrix.") {
        $next = $d.Paragraphs.Item($i + 1)
        $nt = $next.Range.Text.TrimEnd([char]13)
        if ($nt -eq "") {
            $next.Range.Delete()
        }
        break
    }
}

# ---------------------------------------------------------------------
# 3. Merge the two runs "...save the document" + "." into a single run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute(
    "The following Risk Assessment template has been provided in a word format to enable you to type in information and to electronically transmit and save the document.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The following Risk Assessment template has been provided in a word format to enable you to type in information and to electronically transmit and save the document.",
    2
)

# ---------------------------------------------------------------------
# 4. Update the header text "GMIT Civic Engagement" -> "Civic Engagement".
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute(
        "GMIT Civic Engagement", $true, $false, $false, $false, $false,
        $true, 1, $false, "Civic Engagement", 2
    )
}
